$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Date/Time columns as Text so Excel does not auto-convert
# "2026-02-04" / "11:00:00" style strings into date/time serials on write.
$ws.Range("B2:C22").NumberFormat = "@"

$data = New-Object 'object[,]' 21,41

$data[0,0] = 'Algerian Ligue 1'
$data[0,1] = '2026-02-04'
$data[0,2] = '11:00:00'
$data[0,3] = 'USM Khenchela'
$data[0,4] = 'ASO Chlef'
$data[0,5] = 1.8
$data[0,6] = 2.32
$data[0,7] = 4.1
$data[0,8] = 7.6
$data[0,9] = 2.66
$data[0,10] = 5
$data[0,11] = 1.01
$data[0,12] = 1.01
$data[0,13] = 1.97
$data[0,14] = 1.32
$data[0,15] = 1.34
$data[0,16] = 2.46
$data[0,17] = 1.12
$data[0,18] = 1.05
$data[0,19] = 1.05
$data[0,20] = 1.05
$data[0,21] = 1.15
$data[0,22] = 1.75
$data[0,23] = 970
$data[0,24] = 970
$data[0,25] = 1000
$data[0,26] = 1000
$data[0,27] = 970
$data[0,28] = 970
$data[0,29] = 970
$data[0,30] = 1000
$data[0,31] = 1000
$data[0,32] = 970
$data[0,33] = 970
$data[0,34] = 1000
$data[0,35] = 1000
$data[0,36] = 1000
$data[0,37] = 1000
$data[0,38] = 1000
$data[0,39] = 1000
$data[0,40] = 1000

$data[1,0] = 'Algerian Ligue 1'
$data[1,1] = '2026-02-04'
$data[1,2] = '11:00:00'
$data[1,3] = 'MC El Bayadh'
$data[1,4] = 'MB Rouissat'
$data[1,5] = 1.04
$data[1,6] = 1000
$data[1,7] = 1.04
$data[1,8] = 1000
$data[1,9] = 1.02
$data[1,10] = 950
$data[1,11] = 1.01
$data[1,12] = 1.01
$data[1,13] = 1.25
$data[1,14] = 1.01
$data[1,15] = 1.24
$data[1,16] = 1.02
$data[1,17] = 1.12
$data[1,18] = 1.05
$data[1,19] = 1.05
$data[1,20] = 1.05
$data[1,21] = 1.01
$data[1,22] = 1.01
$data[1,23] = 970
$data[1,24] = 970
$data[1,25] = 1000
$data[1,26] = 1000
$data[1,27] = 970
$data[1,28] = 970
$data[1,29] = 970
$data[1,30] = 1000
$data[1,31] = 1000
$data[1,32] = 970
$data[1,33] = 970
$data[1,34] = 1000
$data[1,35] = 1000
$data[1,36] = 1000
$data[1,37] = 1000
$data[1,38] = 1000
$data[1,39] = 1000
$data[1,40] = 1000

$data[2,0] = 'Romanian Liga I'
$data[2,1] = '2026-02-04'
$data[2,2] = '11:00:00'
$data[2,3] = 'FC Metaloglobus Bucuresti'
$data[2,4] = 'Csikszereda'
$data[2,5] = 2.98
$data[2,6] = 3.3
$data[2,7] = 2.4
$data[2,8] = 2.6
$data[2,9] = 3.45
$data[2,10] = 3.75
$data[2,11] = 1.01
$data[2,12] = 1.07
$data[2,13] = 3.45
$data[2,14] = 1.34
$data[2,15] = 1.87
$data[2,16] = 2
$data[2,17] = 1.33
$data[2,18] = 3.5
$data[2,19] = 1.78
$data[2,20] = 2.1
$data[2,21] = 1.62
$data[2,22] = 1.45
$data[2,23] = 970
$data[2,24] = 12
$data[2,25] = 18.5
$data[2,26] = 42
$data[2,27] = 970
$data[2,28] = 9.199999999999999
$data[2,29] = 14
$data[2,30] = 32
$data[2,31] = 25
$data[2,32] = 16
$data[2,33] = 970
$data[2,34] = 48
$data[2,35] = 65
$data[2,36] = 42
$data[2,37] = 55
$data[2,38] = 110
$data[2,39] = 40
$data[2,40] = 26

$data[3,0] = 'Egyptian Premier'
$data[3,1] = '2026-02-04'
$data[3,2] = '12:00:00'
$data[3,3] = 'Smouha'
$data[3,4] = 'Pyramids'
$data[3,5] = 4.2
$data[3,6] = 4.9
$data[3,7] = 2.14
$data[3,8] = 2.22
$data[3,9] = 2.96
$data[3,10] = 3.25
$data[3,11] = 1.01
$data[3,12] = 1.12
$data[3,13] = 2.48
$data[3,14] = 1.55
$data[3,15] = 1.5
$data[3,16] = 2.62
$data[3,17] = 1.17
$data[3,18] = 5.3
$data[3,19] = 2.14
$data[3,20] = 1.7
$data[3,21] = 1.81
$data[3,22] = 1.26
$data[3,23] = 8.6
$data[3,24] = 7.6
$data[3,25] = 12
$data[3,26] = 30
$data[3,27] = 12.5
$data[3,28] = 7.4
$data[3,29] = 12
$data[3,30] = 32
$data[3,31] = 32
$data[3,32] = 21
$data[3,33] = 26
$data[3,34] = 65
$data[3,35] = 120
$data[3,36] = 85
$data[3,37] = 120
$data[3,38] = 240
$data[3,39] = 140
$data[3,40] = 34

$data[4,0] = 'Egyptian Premier'
$data[4,1] = '2026-02-04'
$data[4,2] = '12:00:00'
$data[4,3] = 'ZED FC'
$data[4,4] = 'Al-Masry'
$data[4,5] = 2.92
$data[4,6] = 3.35
$data[4,7] = 2.68
$data[4,8] = 3
$data[4,9] = 2.9
$data[4,10] = 3.3
$data[4,11] = 1.01
$data[4,12] = 1.01
$data[4,13] = 1.91
$data[4,14] = 1.5
$data[4,15] = 1.55
$data[4,16] = 2.54
$data[4,17] = 1.17
$data[4,18] = 5.1
$data[4,19] = 2
$data[4,20] = 1.69
$data[4,21] = 1.5
$data[4,22] = 1.42
$data[4,23] = 10.5
$data[4,24] = 12
$data[4,25] = 1000
$data[4,26] = 55
$data[4,27] = 12.5
$data[4,28] = 9.6
$data[4,29] = 16
$data[4,30] = 48
$data[4,31] = 24
$data[4,32] = 17
$data[4,33] = 26
$data[4,34] = 75
$data[4,35] = 70
$data[4,36] = 55
$data[4,37] = 80
$data[4,38] = 1000
$data[4,39] = 1000
$data[4,40] = 1000

$data[5,0] = 'Romanian Liga I'
$data[5,1] = '2026-02-04'
$data[5,2] = '13:00:00'
$data[5,3] = 'UTA Arad'
$data[5,4] = 'CFR Cluj'
$data[5,5] = 3.05
$data[5,6] = 3.6
$data[5,7] = 2.34
$data[5,8] = 2.66
$data[5,9] = 3.3
$data[5,10] = 3.65
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 1.74
$data[5,16] = 2.08
$data[5,17] = 0
$data[5,18] = 0
$data[5,19] = 0
$data[5,20] = 0
$data[5,21] = 0
$data[5,22] = 0
$data[5,23] = 0
$data[5,24] = 0
$data[5,25] = 0
$data[5,26] = 0
$data[5,27] = 0
$data[5,28] = 0
$data[5,29] = 0
$data[5,30] = 0
$data[5,31] = 0
$data[5,32] = 0
$data[5,33] = 0
$data[5,34] = 0
$data[5,35] = 0
$data[5,36] = 0
$data[5,37] = 0
$data[5,38] = 0
$data[5,39] = 0
$data[5,40] = 0

$data[6,0] = 'Slovenian Premier League'
$data[6,1] = '2026-02-04'
$data[6,2] = '13:30:00'
$data[6,3] = 'NK Celje'
$data[6,4] = 'NK Radomlje'
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 1.24
$data[6,16] = 1.01
$data[6,17] = 0
$data[6,18] = 0
$data[6,19] = 0
$data[6,20] = 0
$data[6,21] = 0
$data[6,22] = 0
$data[6,23] = 0
$data[6,24] = 0
$data[6,25] = 0
$data[6,26] = 0
$data[6,27] = 0
$data[6,28] = 0
$data[6,29] = 0
$data[6,30] = 0
$data[6,31] = 0
$data[6,32] = 0
$data[6,33] = 0
$data[6,34] = 0
$data[6,35] = 0
$data[6,36] = 0
$data[6,37] = 0
$data[6,38] = 0
$data[6,39] = 0
$data[6,40] = 0

$data[7,0] = 'Egyptian Premier'
$data[7,1] = '2026-02-04'
$data[7,2] = '15:00:00'
$data[7,3] = 'Kahraba Ismailia'
$data[7,4] = 'Zamalek'
$data[7,5] = 5.7
$data[7,6] = 6.8
$data[7,7] = 1.7
$data[7,8] = 1.81
$data[7,9] = 3.45
$data[7,10] = 3.95
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 0
$data[7,14] = 0
$data[7,15] = 1.67
$data[7,16] = 2.3
$data[7,17] = 0
$data[7,18] = 0
$data[7,19] = 0
$data[7,20] = 0
$data[7,21] = 0
$data[7,22] = 0
$data[7,23] = 0
$data[7,24] = 0
$data[7,25] = 0
$data[7,26] = 0
$data[7,27] = 0
$data[7,28] = 0
$data[7,29] = 0
$data[7,30] = 0
$data[7,31] = 0
$data[7,32] = 0
$data[7,33] = 0
$data[7,34] = 0
$data[7,35] = 0
$data[7,36] = 0
$data[7,37] = 0
$data[7,38] = 0
$data[7,39] = 0
$data[7,40] = 0

$data[8,0] = 'Romanian Liga I'
$data[8,1] = '2026-02-04'
$data[8,2] = '15:30:00'
$data[8,3] = 'Farul Constanta'
$data[8,4] = 'Dinamo Bucharest'
$data[8,5] = 3.75
$data[8,6] = 4.2
$data[8,7] = 2.04
$data[8,8] = 2.18
$data[8,9] = 3.45
$data[8,10] = 3.8
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 1.92
$data[8,16] = 1.95
$data[8,17] = 0
$data[8,18] = 0
$data[8,19] = 0
$data[8,20] = 0
$data[8,21] = 0
$data[8,22] = 0
$data[8,23] = 0
$data[8,24] = 0
$data[8,25] = 0
$data[8,26] = 0
$data[8,27] = 0
$data[8,28] = 0
$data[8,29] = 0
$data[8,30] = 0
$data[8,31] = 0
$data[8,32] = 0
$data[8,33] = 0
$data[8,34] = 0
$data[8,35] = 0
$data[8,36] = 0
$data[8,37] = 0
$data[8,38] = 0
$data[8,39] = 0
$data[8,40] = 0

$data[9,0] = 'Scottish Premiership'
$data[9,1] = '2026-02-04'
$data[9,2] = '16:45:00'
$data[9,3] = 'Livingston'
$data[9,4] = 'Falkirk'
$data[9,5] = 3.5
$data[9,6] = 3.7
$data[9,7] = 2.22
$data[9,8] = 2.26
$data[9,9] = 3.55
$data[9,10] = 3.65
$data[9,11] = 0
$data[9,12] = 1.07
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 1.86
$data[9,16] = 2.02
$data[9,17] = 0
$data[9,18] = 0
$data[9,19] = 1.78
$data[9,20] = 2.08
$data[9,21] = 0
$data[9,22] = 0
$data[9,23] = 14
$data[9,24] = 10
$data[9,25] = 17.5
$data[9,26] = 32
$data[9,27] = 14
$data[9,28] = 8.4
$data[9,29] = 13.5
$data[9,30] = 32
$data[9,31] = 27
$data[9,32] = 16
$data[9,33] = 19
$data[9,34] = 1000
$data[9,35] = 1000
$data[9,36] = 1000
$data[9,37] = 60
$data[9,38] = 130
$data[9,39] = 1000
$data[9,40] = 20

$data[10,0] = 'Scottish Premiership'
$data[10,1] = '2026-02-04'
$data[10,2] = '16:45:00'
$data[10,3] = 'Dundee'
$data[10,4] = 'Motherwell'
$data[10,5] = 4.6
$data[10,6] = 4.7
$data[10,7] = 1.94
$data[10,8] = 2
$data[10,9] = 3.45
$data[10,10] = 3.75
$data[10,11] = 0
$data[10,12] = 1.08
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 1.86
$data[10,16] = 1.98
$data[10,17] = 0
$data[10,18] = 0
$data[10,19] = 1.84
$data[10,20] = 2.02
$data[10,21] = 0
$data[10,22] = 0
$data[10,23] = 16
$data[10,24] = 1000
$data[10,25] = 13
$data[10,26] = 26
$data[10,27] = 1000
$data[10,28] = 9.199999999999999
$data[10,29] = 1000
$data[10,30] = 980
$data[10,31] = 980
$data[10,32] = 23
$data[10,33] = 22
$data[10,34] = 42
$data[10,35] = 120
$data[10,36] = 65
$data[10,37] = 85
$data[10,38] = 120
$data[10,39] = 75
$data[10,40] = 16

$data[11,0] = 'Scottish Premiership'
$data[11,1] = '2026-02-04'
$data[11,2] = '16:45:00'
$data[11,3] = 'Rangers'
$data[11,4] = 'Kilmarnock'
$data[11,5] = 1.32
$data[11,6] = 1.37
$data[11,7] = 11
$data[11,8] = 12.5
$data[11,9] = 5.6
$data[11,10] = 6.2
$data[11,11] = 0
$data[11,12] = 1.04
$data[11,13] = 5.2
$data[11,14] = 1.22
$data[11,15] = 2.46
$data[11,16] = 1.6
$data[11,17] = 1.58
$data[11,18] = 2.52
$data[11,19] = 2.02
$data[11,20] = 1.85
$data[11,21] = 0
$data[11,22] = 0
$data[11,23] = 26
$data[11,24] = 1000
$data[11,25] = 130
$data[11,26] = 470
$data[11,27] = 9.800000000000001
$data[11,28] = 13.5
$data[11,29] = 1000
$data[11,30] = 200
$data[11,31] = 8.800000000000001
$data[11,32] = 11
$data[11,33] = 1000
$data[11,34] = 150
$data[11,35] = 11
$data[11,36] = 1000
$data[11,37] = 1000
$data[11,38] = 160
$data[11,39] = 4.9
$data[11,40] = 1000

$data[12,0] = 'Scottish Premiership'
$data[12,1] = '2026-02-04'
$data[12,2] = '16:45:00'
$data[12,3] = 'Hibernian'
$data[12,4] = 'Dundee Utd'
$data[12,5] = 1.69
$data[12,6] = 1.7
$data[12,7] = 5.6
$data[12,8] = 6.4
$data[12,9] = 3.95
$data[12,10] = 4.3
$data[12,11] = 0
$data[12,12] = 1.06
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 2.04
$data[12,16] = 1.84
$data[12,17] = 0
$data[12,18] = 0
$data[12,19] = 1.84
$data[12,20] = 2.04
$data[12,21] = 0
$data[12,22] = 0
$data[12,23] = 17.5
$data[12,24] = 1000
$data[12,25] = 1000
$data[12,26] = 170
$data[12,27] = 8.800000000000001
$data[12,28] = 9.199999999999999
$data[12,29] = 1000
$data[12,30] = 1000
$data[12,31] = 11
$data[12,32] = 9.800000000000001
$data[12,33] = 1000
$data[12,34] = 1000
$data[12,35] = 1000
$data[12,36] = 1000
$data[12,37] = 1000
$data[12,38] = 130
$data[12,39] = 9.800000000000001
$data[12,40] = 1000

$data[13,0] = 'Colombian Primera B'
$data[13,1] = '2026-02-04'
$data[13,2] = '17:00:00'
$data[13,3] = 'Tigres FC Zipaquira'
$data[13,4] = 'Orsomarso'
$data[13,5] = 1.04
$data[13,6] = 1000
$data[13,7] = 1.04
$data[13,8] = 1000
$data[13,9] = 1.01
$data[13,10] = 1000
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 0
$data[13,14] = 0
$data[13,15] = 1.24
$data[13,16] = 1.01
$data[13,17] = 0
$data[13,18] = 0
$data[13,19] = 0
$data[13,20] = 0
$data[13,21] = 0
$data[13,22] = 0
$data[13,23] = 0
$data[13,24] = 0
$data[13,25] = 0
$data[13,26] = 0
$data[13,27] = 0
$data[13,28] = 0
$data[13,29] = 0
$data[13,30] = 0
$data[13,31] = 0
$data[13,32] = 0
$data[13,33] = 0
$data[13,34] = 0
$data[13,35] = 0
$data[13,36] = 0
$data[13,37] = 0
$data[13,38] = 0
$data[13,39] = 0
$data[13,40] = 0

$data[14,0] = 'Scottish Premiership'
$data[14,1] = '2026-02-04'
$data[14,2] = '17:00:00'
$data[14,3] = 'Aberdeen'
$data[14,4] = 'Celtic'
$data[14,5] = 6
$data[14,6] = 6.8
$data[14,7] = 1.55
$data[14,8] = 1.6
$data[14,9] = 4.5
$data[14,10] = 5
$data[14,11] = 0
$data[14,12] = 1.04
$data[14,13] = 5
$data[14,14] = 1.22
$data[14,15] = 2.36
$data[14,16] = 1.66
$data[14,17] = 1.54
$data[14,18] = 2.62
$data[14,19] = 1.78
$data[14,20] = 2.14
$data[14,21] = 0
$data[14,22] = 0
$data[14,23] = 28
$data[14,24] = 10.5
$data[14,25] = 10.5
$data[14,26] = 15.5
$data[14,27] = 1000
$data[14,28] = 11
$data[14,29] = 10.5
$data[14,30] = 1000
$data[14,31] = 1000
$data[14,32] = 27
$data[14,33] = 1000
$data[14,34] = 1000
$data[14,35] = 200
$data[14,36] = 1000
$data[14,37] = 1000
$data[14,38] = 1000
$data[14,39] = 1000
$data[14,40] = 6.8

$data[15,0] = 'Brazilian Serie A'
$data[15,1] = '2026-02-04'
$data[15,2] = '19:00:00'
$data[15,3] = 'Red Bull Bragantino'
$data[15,4] = 'Atletico MG'
$data[15,5] = 2.26
$data[15,6] = 2.38
$data[15,7] = 3.5
$data[15,8] = 3.8
$data[15,9] = 3.3
$data[15,10] = 3.4
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 0
$data[15,14] = 0
$data[15,15] = 1.78
$data[15,16] = 2.24
$data[15,17] = 0
$data[15,18] = 0
$data[15,19] = 0
$data[15,20] = 0
$data[15,21] = 0
$data[15,22] = 0
$data[15,23] = 0
$data[15,24] = 0
$data[15,25] = 0
$data[15,26] = 0
$data[15,27] = 0
$data[15,28] = 0
$data[15,29] = 0
$data[15,30] = 0
$data[15,31] = 0
$data[15,32] = 0
$data[15,33] = 0
$data[15,34] = 0
$data[15,35] = 0
$data[15,36] = 0
$data[15,37] = 0
$data[15,38] = 0
$data[15,39] = 0
$data[15,40] = 0

$data[16,0] = 'Brazilian Serie A'
$data[16,1] = '2026-02-04'
$data[16,2] = '19:00:00'
$data[16,3] = 'Flamengo'
$data[16,4] = 'Internacional'
$data[16,5] = 1.44
$data[16,6] = 1.48
$data[16,7] = 9
$data[16,8] = 11
$data[16,9] = 4.6
$data[16,10] = 5.1
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 0
$data[16,14] = 0
$data[16,15] = 2.04
$data[16,16] = 1.89
$data[16,17] = 0
$data[16,18] = 0
$data[16,19] = 0
$data[16,20] = 0
$data[16,21] = 0
$data[16,22] = 0
$data[16,23] = 0
$data[16,24] = 0
$data[16,25] = 0
$data[16,26] = 0
$data[16,27] = 0
$data[16,28] = 0
$data[16,29] = 0
$data[16,30] = 0
$data[16,31] = 0
$data[16,32] = 0
$data[16,33] = 0
$data[16,34] = 0
$data[16,35] = 0
$data[16,36] = 0
$data[16,37] = 0
$data[16,38] = 0
$data[16,39] = 0
$data[16,40] = 0

$data[17,0] = 'Colombian Primera B'
$data[17,1] = '2026-02-04'
$data[17,2] = '19:30:00'
$data[17,3] = 'Internacional de Palmira'
$data[17,4] = 'Quindio'
$data[17,5] = 1.04
$data[17,6] = 1000
$data[17,7] = 1.04
$data[17,8] = 1000
$data[17,9] = 1.01
$data[17,10] = 1000
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 0
$data[17,14] = 0
$data[17,15] = 1.24
$data[17,16] = 1.01
$data[17,17] = 0
$data[17,18] = 0
$data[17,19] = 0
$data[17,20] = 0
$data[17,21] = 0
$data[17,22] = 0
$data[17,23] = 0
$data[17,24] = 0
$data[17,25] = 0
$data[17,26] = 0
$data[17,27] = 0
$data[17,28] = 0
$data[17,29] = 0
$data[17,30] = 0
$data[17,31] = 0
$data[17,32] = 0
$data[17,33] = 0
$data[17,34] = 0
$data[17,35] = 0
$data[17,36] = 0
$data[17,37] = 0
$data[17,38] = 0
$data[17,39] = 0
$data[17,40] = 0

$data[18,0] = 'Brazilian Serie A'
$data[18,1] = '2026-02-04'
$data[18,2] = '20:00:00'
$data[18,3] = 'Remo'
$data[18,4] = 'Mirassol'
$data[18,5] = 2.96
$data[18,6] = 3.2
$data[18,7] = 2.62
$data[18,8] = 2.8
$data[18,9] = 3.25
$data[18,10] = 3.3
$data[18,11] = 0
$data[18,12] = 0
$data[18,13] = 0
$data[18,14] = 0
$data[18,15] = 1.68
$data[18,16] = 2.3
$data[18,17] = 0
$data[18,18] = 0
$data[18,19] = 0
$data[18,20] = 0
$data[18,21] = 0
$data[18,22] = 0
$data[18,23] = 0
$data[18,24] = 0
$data[18,25] = 0
$data[18,26] = 0
$data[18,27] = 0
$data[18,28] = 0
$data[18,29] = 0
$data[18,30] = 0
$data[18,31] = 0
$data[18,32] = 0
$data[18,33] = 0
$data[18,34] = 0
$data[18,35] = 0
$data[18,36] = 0
$data[18,37] = 0
$data[18,38] = 0
$data[18,39] = 0
$data[18,40] = 0

$data[19,0] = 'Brazilian Serie A'
$data[19,1] = '2026-02-04'
$data[19,2] = '20:00:00'
$data[19,3] = 'Santos'
$data[19,4] = 'Sao Paulo'
$data[19,5] = 2.18
$data[19,6] = 2.28
$data[19,7] = 3.7
$data[19,8] = 4.1
$data[19,9] = 3.3
$data[19,10] = 3.45
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 0
$data[19,14] = 0
$data[19,15] = 1.71
$data[19,16] = 2.28
$data[19,17] = 0
$data[19,18] = 0
$data[19,19] = 0
$data[19,20] = 0
$data[19,21] = 0
$data[19,22] = 0
$data[19,23] = 0
$data[19,24] = 0
$data[19,25] = 0
$data[19,26] = 0
$data[19,27] = 0
$data[19,28] = 0
$data[19,29] = 0
$data[19,30] = 0
$data[19,31] = 0
$data[19,32] = 0
$data[19,33] = 0
$data[19,34] = 0
$data[19,35] = 0
$data[19,36] = 0
$data[19,37] = 0
$data[19,38] = 0
$data[19,39] = 0
$data[19,40] = 0

$data[20,0] = 'Colombian Primera B'
$data[20,1] = '2026-02-04'
$data[20,2] = '20:00:00'
$data[20,3] = 'Independiente Yumbo'
$data[20,4] = 'Real Cartagena'
$data[20,5] = 1.04
$data[20,6] = 1000
$data[20,7] = 1.04
$data[20,8] = 1000
$data[20,9] = 1.01
$data[20,10] = 1000
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 0
$data[20,14] = 0
$data[20,15] = 1.24
$data[20,16] = 1.01
$data[20,17] = 0
$data[20,18] = 0
$data[20,19] = 0
$data[20,20] = 0
$data[20,21] = 0
$data[20,22] = 0
$data[20,23] = 0
$data[20,24] = 0
$data[20,25] = 0
$data[20,26] = 0
$data[20,27] = 0
$data[20,28] = 0
$data[20,29] = 0
$data[20,30] = 0
$data[20,31] = 0
$data[20,32] = 0
$data[20,33] = 0
$data[20,34] = 0
$data[20,35] = 0
$data[20,36] = 0
$data[20,37] = 0
$data[20,38] = 0
$data[20,39] = 0
$data[20,40] = 0

$ws.Range("A2:AO22").Value = $data

# Restore default (unstyled) formatting on the Date/Time columns so the
# written cells match the rest of the sheet (no explicit cell style).
$ws.Range("B2:C22").Style = "Normal"